# Update the "Metadata" sheet with the new URL, Version, Date, and Publisher
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-use-scale"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# On the "Elements" sheet, clear the stray "Constraint(s)" value that had been
# duplicated onto the root Extension row (row 2) - it belongs only on the
# Extension.extension row (row 4).
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
